# refactor del calculo de riesgo y pureza proporcion
#
# This particular re-run of the pcsmote logging pass only changed the
# capture "timestamp" column (Z) for every data row (rows 2-97) on the
# active/only worksheet. All other columns (A-Y) are untouched.
#
# The new timestamps were captured in tight batches during the run, so
# many consecutive rows share the exact same instant - those are grouped
# into a single range write below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z37").Value  = "2025-11-05T14:13:22.032933"
$ws.Range("Z38:Z42").Value = "2025-11-05T14:13:22.042064"
$ws.Range("Z43:Z54").Value = "2025-11-05T14:13:22.042587"
$ws.Range("Z55").Value     = "2025-11-05T14:13:22.221975"
$ws.Range("Z56").Value     = "2025-11-05T14:13:22.222887"
$ws.Range("Z57:Z59").Value = "2025-11-05T14:13:22.223887"
$ws.Range("Z60:Z63").Value = "2025-11-05T14:13:22.224900"
$ws.Range("Z64:Z65").Value = "2025-11-05T14:13:22.225976"
$ws.Range("Z66:Z67").Value = "2025-11-05T14:13:22.226505"
$ws.Range("Z68:Z71").Value = "2025-11-05T14:13:22.491398"
$ws.Range("Z72:Z76").Value = "2025-11-05T14:13:22.492394"
$ws.Range("Z77:Z84").Value = "2025-11-05T14:13:22.733741"
$ws.Range("Z85:Z95").Value = "2025-11-05T14:13:22.734742"
$ws.Range("Z96:Z97").Value = "2025-11-05T14:13:22.735741"
